$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as TEXT (avoids Excel auto-converting
# numeric-looking strings like "1.000" or "0.9999" into real numbers), while
# leaving the cell style untouched (ClearFormats drops the temporary "@" text
# format we apply so the cell ends up with no explicit style, same as before).
function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

Set-TextValue $ws.Range('D2') '26.942.40'
Set-TextValue $ws.Range('E2') '  -3.24%  '
Set-TextValue $ws.Range('D3') '1.739.31'
Set-TextValue $ws.Range('E3') '  -1.35%  '
Set-TextValue $ws.Range('E4') '  -0.03%  '
Set-TextValue $ws.Range('D5') '310.67'
Set-TextValue $ws.Range('E5') '  -5.57%  '
Set-TextValue $ws.Range('D6') '1.000'
Set-TextValue $ws.Range('E6') '  -0.08%  '
Set-TextValue $ws.Range('D7') '0.4975'
Set-TextValue $ws.Range('E7') '  +3.14%  '
Set-TextValue $ws.Range('D8') '0.3557'
Set-TextValue $ws.Range('E8') '  -0.04%  '
Set-TextValue $ws.Range('D9') '42.51'
Set-TextValue $ws.Range('E9') '  -1.60%  '
Set-TextValue $ws.Range('D10') '0.07252'
Set-TextValue $ws.Range('E10') '  -4.10%  '
Set-TextValue $ws.Range('D11') '1.059'
Set-TextValue $ws.Range('E11') '  -2.29%  '
Set-TextValue $ws.Range('D12') '1.000'
Set-TextValue $ws.Range('E12') '  -0.03%  '
Set-TextValue $ws.Range('D13') '19.96'
Set-TextValue $ws.Range('E13') '  -3.01%  '
Set-TextValue $ws.Range('D14') '5.964'
Set-TextValue $ws.Range('E14') '  -2.28%  '
Set-TextValue $ws.Range('D15') '1.733.68'
Set-TextValue $ws.Range('E15') '  -1.60%  '
Set-TextValue $ws.Range('D16') '6.846'
Set-TextValue $ws.Range('E16') '  -4.34%  '
Set-TextValue $ws.Range('D17') '86.46'
Set-TextValue $ws.Range('E17') '  -7.31%  '
Set-TextValue $ws.Range('D18') '0.00001034'
Set-TextValue $ws.Range('E18') '  -5.18%  '
Set-TextValue $ws.Range('D19') '0.06387'
Set-TextValue $ws.Range('E19') '  -0.70%  '
Set-TextValue $ws.Range('E20') '  -0.07%  '
Set-TextValue $ws.Range('D21') '16.58'
Set-TextValue $ws.Range('E21') '  -1.51%  '
Set-TextValue $ws.Range('D22') '5.737'
Set-TextValue $ws.Range('E22') '  -1.31%  '
Set-TextValue $ws.Range('D23') '27.004.99'
Set-TextValue $ws.Range('E23') '  -3.09%  '
Set-TextValue $ws.Range('D24') '11.20'
Set-TextValue $ws.Range('E24') '  +0.40%  '
Set-TextValue $ws.Range('D25') '2.045'
Set-TextValue $ws.Range('E25') '  -5.40%  '
Set-TextValue $ws.Range('D26') '153.49'
Set-TextValue $ws.Range('E26') '  -6.46%  '
Set-TextValue $ws.Range('D27') '19.88'
Set-TextValue $ws.Range('E27') '  -1.39%  '
Set-TextValue $ws.Range('D28') '1.934.49'
Set-TextValue $ws.Range('E28') '  -1.55%  '
Set-TextValue $ws.Range('D29') '2.131'
Set-TextValue $ws.Range('E29') '  -3.51%  '
Set-TextValue $ws.Range('D30') '120.69'
Set-TextValue $ws.Range('E30') '  -2.19%  '
Set-TextValue $ws.Range('D31') '1.058'
Set-TextValue $ws.Range('E31') '  -0.08%  '
Set-TextValue $ws.Range('D32') '0.09414'
Set-TextValue $ws.Range('E32') '  -0.62%  '
Set-TextValue $ws.Range('D33') '3.571'
Set-TextValue $ws.Range('E33') '  -2.25%  '
Set-TextValue $ws.Range('D34') '5.378'
Set-TextValue $ws.Range('E34') '  -3.34%  '
Set-TextValue $ws.Range('D35') '0.02196'
Set-TextValue $ws.Range('E35') '  -3.27%  '
Set-TextValue $ws.Range('D36') '0.05904'
Set-TextValue $ws.Range('E36') '  -1.55%  '
Set-TextValue $ws.Range('D37') '11.05'
Set-TextValue $ws.Range('E37') '  -5.13%  '
Set-TextValue $ws.Range('D38') '0.1999'
Set-TextValue $ws.Range('E38') '  -3.46%  '
Set-TextValue $ws.Range('D39') '1.423'
Set-TextValue $ws.Range('E39') '  -0.63%  '
Set-TextValue $ws.Range('D40') '4.760'
Set-TextValue $ws.Range('E40') '  -2.80%  '
Set-TextValue $ws.Range('D41') '0.9999'
Set-TextValue $ws.Range('E41') '  -0.11%  '
Set-TextValue $ws.Range('D42') '0.6009'
Set-TextValue $ws.Range('E42') '  -2.56%  '
Set-TextValue $ws.Range('D43') '1.108'
Set-TextValue $ws.Range('E43') '  -6.36%  '
Set-TextValue $ws.Range('D44') '7.460'
Set-TextValue $ws.Range('E44') '  -3.76%  '
Set-TextValue $ws.Range('D45') '12.87'
Set-TextValue $ws.Range('E45') '  -1.40%  '
Set-TextValue $ws.Range('D46') '3.575'
Set-TextValue $ws.Range('E46') '  -4.36%  '
Set-TextValue $ws.Range('D47') '0.5633'
Set-TextValue $ws.Range('E47') '  -3.04%  '
Set-TextValue $ws.Range('D48') '119.89'
Set-TextValue $ws.Range('E48') '  -2.93%  '
Set-TextValue $ws.Range('D49') '1.855'
Set-TextValue $ws.Range('E49') '  -3.68%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue $ws.Range('D50') '1.100'
Set-TextValue $ws.Range('E50') '  -4.43%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D51') '0.06662'
Set-TextValue $ws.Range('E51') '  -1.92%  '
